$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.025109979222057
$ws.Cells.Item(2, 4).Value = 1.028327507728955
$ws.Cells.Item(2, 5).Value = 1.025421719044825
$ws.Cells.Item(2, 6).Value = 1.023603807099087
$ws.Cells.Item(2, 9).Value = 1.029356175659613
$ws.Cells.Item(2, 10).Value = 1.030280834624214
$ws.Cells.Item(2, 11).Value = 1.031144646376276
$ws.Cells.Item(2, 12).Value = 1.028247338075449
$ws.Cells.Item(2, 13).Value = 1.026434758240348
$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.026816996821466
$ws.Cells.Item(3, 4).Value = 1.029572503769889
$ws.Cells.Item(3, 5).Value = 1.026897240026056
$ws.Cells.Item(3, 6).Value = 1.025958908886571
$ws.Cells.Item(3, 9).Value = 1.029763678090753
$ws.Cells.Item(3, 10).Value = 1.031623137731457
$ws.Cells.Item(3, 11).Value = 1.032196375637372
$ws.Cells.Item(3, 12).Value = 1.029528345412348
$ws.Cells.Item(3, 13).Value = 1.02859256105393
$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.027917597803258
$ws.Cells.Item(4, 4).Value = 1.030374573288738
$ws.Cells.Item(4, 5).Value = 1.027848682398705
$ws.Cells.Item(4, 6).Value = 1.027478105170414
$ws.Cells.Item(4, 9).Value = 1.030024036439695
$ws.Cells.Item(4, 10).Value = 1.0324875167029
$ws.Cells.Item(4, 11).Value = 1.032872847853746
$ws.Cells.Item(4, 12).Value = 1.030353437239987
$ws.Cells.Item(4, 13).Value = 1.029983813633798
$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.028379362758668
$ws.Cells.Item(5, 4).Value = 1.030710931655675
$ws.Cells.Item(5, 5).Value = 1.028247889577853
$ws.Cells.Item(5, 6).Value = 1.028115678875713
$ws.Cells.Item(5, 9).Value = 1.030132700974056
$ws.Cells.Item(5, 10).Value = 1.032849916000458
$ws.Cells.Item(5, 11).Value = 1.033156274971455
$ws.Cells.Item(5, 12).Value = 1.030699409022867
$ws.Cells.Item(5, 13).Value = 1.030567530743169
$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.028456841267622
$ws.Cells.Item(6, 4).Value = 1.030767359297107
$ws.Cells.Item(6, 5).Value = 1.028314872966502
$ws.Cells.Item(6, 6).Value = 1.02822266692482
$ws.Cells.Item(6, 9).Value = 1.030150900050432
$ws.Cells.Item(6, 10).Value = 1.032910707080482
$ws.Cells.Item(6, 11).Value = 1.033203807571908
$ws.Cells.Item(6, 12).Value = 1.030757447027979
$ws.Cells.Item(6, 13).Value = 1.030665471911107
$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.027923771548911
$ws.Cells.Item(7, 4).Value = 1.030379070976491
$ws.Cells.Item(7, 5).Value = 1.027854019665564
$ws.Cells.Item(7, 6).Value = 1.027486628726976
$ws.Cells.Item(7, 9).Value = 1.030025491517732
$ws.Cells.Item(7, 10).Value = 1.032492362952099
$ws.Cells.Item(7, 11).Value = 1.032876638783877
$ws.Cells.Item(7, 12).Value = 1.030358063636423
$ws.Cells.Item(7, 13).Value = 1.029991617830588
$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.025687703796404
$ws.Cells.Item(8, 4).Value = 1.028748997667068
$ws.Cells.Item(8, 5).Value = 1.025921074579339
$ws.Cells.Item(8, 6).Value = 1.024400719605939
$ws.Cells.Item(8, 9).Value = 1.029494584338202
$ws.Cells.Item(8, 10).Value = 1.030735347187083
$ws.Cells.Item(8, 11).Value = 1.031500932459475
$ws.Cells.Item(8, 12).Value = 1.028681057095857
$ws.Cells.Item(8, 13).Value = 1.027165048704868
$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.021716337316919
$ws.Cells.Item(9, 4).Value = 1.025849027034356
$ws.Cells.Item(9, 5).Value = 1.022488874597086
$ws.Cells.Item(9, 6).Value = 1.018925432794981
$ws.Cells.Item(9, 9).Value = 1.028533362337054
$ws.Cells.Item(9, 10).Value = 1.027606548005155
$ws.Cells.Item(9, 11).Value = 1.029045078098312
$ws.Cells.Item(9, 12).Value = 1.025696166649526
$ws.Cells.Item(9, 13).Value = 1.02214473490405
$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.019046632427948
$ws.Cells.Item(10, 4).Value = 1.023896373917577
$ws.Cells.Item(10, 5).Value = 1.020182219179631
$ws.Cells.Item(10, 6).Value = 1.015247998947158
$ws.Cells.Item(10, 9).Value = 1.027874923260944
$ws.Cells.Item(10, 10).Value = 1.025497711254171
$ws.Cells.Item(10, 11).Value = 1.027385782129007
$ws.Cells.Item(10, 12).Value = 1.023685301407006
$ws.Cells.Item(10, 13).Value = 1.01876941823197
$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.0178851033332
$ws.Cells.Item(11, 4).Value = 1.023046086031799
$ws.Cells.Item(11, 5).Value = 1.019178800278517
$ws.Cells.Item(11, 6).Value = 1.013648715430646
$ws.Cells.Item(11, 9).Value = 1.027585552823296
$ws.Cells.Item(11, 10).Value = 1.024578893659922
$ws.Cells.Item(11, 11).Value = 1.026661885708661
$ws.Cells.Item(11, 12).Value = 1.022809403873174
$ws.Cells.Item(11, 13).Value = 1.017300705986151
$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.017452805628249
$ws.Cells.Item(12, 4).Value = 1.022729517949614
$ws.Cells.Item(12, 5).Value = 1.018805372446812
$ws.Cells.Item(12, 6).Value = 1.013053589902746
$ws.Cells.Item(12, 9).Value = 1.027477420534274
$ws.Cells.Item(12, 10).Value = 1.024236731824292
$ws.Cells.Item(12, 11).Value = 1.026392170625026
$ws.Cells.Item(12, 12).Value = 1.022483260572662
$ws.Cells.Item(12, 13).Value = 1.016754046313995
$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.017545573981414
$ws.Cells.Item(13, 4).Value = 1.02279745631656
$ws.Cells.Item(13, 5).Value = 1.018885506573952
$ws.Cells.Item(13, 6).Value = 1.013181295876004
$ws.Cells.Item(13, 9).Value = 1.027500644677788
$ws.Cells.Item(13, 10).Value = 1.024310166472512
$ws.Cells.Item(13, 11).Value = 1.026450063123629
$ws.Cells.Item(13, 12).Value = 1.022553255724774
$ws.Cells.Item(13, 13).Value = 1.016871357691545
$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.017849387029558
$ws.Cells.Item(14, 4).Value = 1.023019933466438
$ws.Cells.Item(14, 5).Value = 1.019147947280976
$ws.Cells.Item(14, 6).Value = 1.013599544438194
$ws.Cells.Item(14, 9).Value = 1.027576627818142
$ws.Cells.Item(14, 10).Value = 1.024550628357462
$ws.Cells.Item(14, 11).Value = 1.026639607947942
$ws.Cells.Item(14, 12).Value = 1.022782461107331
$ws.Cells.Item(14, 13).Value = 1.017255541874831
$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.018036462396112
$ws.Cells.Item(15, 4).Value = 1.023156911336475
$ws.Cells.Item(15, 5).Value = 1.01930955056009
$ws.Cells.Item(15, 6).Value = 1.013857096831287
$ws.Cells.Item(15, 9).Value = 1.027623357566377
$ws.Cells.Item(15, 10).Value = 1.024698668726322
$ws.Cells.Item(15, 11).Value = 1.026756282639454
$ws.Cells.Item(15, 12).Value = 1.022923576099638
$ws.Cells.Item(15, 13).Value = 1.017492101649632
$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.019123600806407
$ws.Cells.Item(16, 4).Value = 1.023952702758047
$ws.Cells.Item(16, 5).Value = 1.020248713795414
$ws.Cells.Item(16, 6).Value = 1.015353988778982
$ws.Cells.Item(16, 9).Value = 1.027894037441805
$ws.Cells.Item(16, 10).Value = 1.025558568911837
$ws.Cells.Item(16, 11).Value = 1.027433709496295
$ws.Cells.Item(16, 12).Value = 1.023743321201805
$ws.Cells.Item(16, 13).Value = 1.018866737524127
$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.019804037279295
$ws.Cells.Item(17, 4).Value = 1.024450591940156
$ws.Cells.Item(17, 5).Value = 1.02083657535601
$ws.Cells.Item(17, 6).Value = 1.016291066600944
$ws.Cells.Item(17, 9).Value = 1.028062681883372
$ws.Cells.Item(17, 10).Value = 1.02609642818957
$ws.Cells.Item(17, 11).Value = 1.027857182538998
$ws.Cells.Item(17, 12).Value = 1.024256126451415
$ws.Cells.Item(17, 13).Value = 1.019727063227277
$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.020200392464898
$ws.Cells.Item(18, 4).Value = 1.024740542518065
$ws.Cells.Item(18, 5).Value = 1.021179020410994
$ws.Cells.Item(18, 6).Value = 1.016836982236367
$ws.Cells.Item(18, 9).Value = 1.028160638587447
$ws.Cells.Item(18, 10).Value = 1.026409605928965
$ws.Cells.Item(18, 11).Value = 1.02810366620078
$ws.Cells.Item(18, 12).Value = 1.024554738534573
$ws.Cells.Item(18, 13).Value = 1.020228186565726
$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.020335449823456
$ws.Cells.Item(19, 4).Value = 1.024839330661269
$ws.Cells.Item(19, 5).Value = 1.021295710475441
$ws.Cells.Item(19, 6).Value = 1.017023013571992
$ws.Cells.Item(19, 9).Value = 1.02819396982979
$ws.Cells.Item(19, 10).Value = 1.026516299326577
$ws.Cells.Item(19, 11).Value = 1.028187622835848
$ws.Cells.Item(19, 12).Value = 1.024656473510397
$ws.Cells.Item(19, 13).Value = 1.020398940708821
$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.019731088070227
$ws.Cells.Item(20, 4).Value = 1.024397220811134
$ws.Cells.Item(20, 5).Value = 1.020773549502739
$ws.Cells.Item(20, 6).Value = 1.016190596189034
$ws.Cells.Item(20, 9).Value = 1.028044630465902
$ws.Cells.Item(20, 10).Value = 1.026038777641628
$ws.Cells.Item(20, 11).Value = 1.027811801885784
$ws.Cells.Item(20, 12).Value = 1.024201158995921
$ws.Cells.Item(20, 13).Value = 1.0196348300606
$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.017759945425523
$ws.Cells.Item(21, 4).Value = 1.022954439861219
$ws.Cells.Item(21, 5).Value = 1.019070684880479
$ws.Cells.Item(21, 6).Value = 1.013476410804059
$ws.Cells.Item(21, 9).Value = 1.027554270599782
$ws.Cells.Item(21, 10).Value = 1.024479842554068
$ws.Cells.Item(21, 11).Value = 1.026583814691333
$ws.Cells.Item(21, 12).Value = 1.0227149879514
$ws.Cells.Item(21, 13).Value = 1.017142440268609
$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.016515656519328
$ws.Cells.Item(22, 4).Value = 1.022043055576406
$ws.Cells.Item(22, 5).Value = 1.01799588932371
$ws.Cells.Item(22, 6).Value = 1.011763623867019
$ws.Cells.Item(22, 9).Value = 1.027242213739954
$ws.Cells.Item(22, 10).Value = 1.023494622974745
$ws.Cells.Item(22, 11).Value = 1.02580693391507
$ws.Cells.Item(22, 12).Value = 1.021775959062057
$ws.Cells.Item(22, 13).Value = 1.015568909117158
$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.017175755021655
$ws.Cells.Item(23, 4).Value = 1.022526605859108
$ws.Cells.Item(23, 5).Value = 1.018566057337284
$ws.Cells.Item(23, 6).Value = 1.012672212286494
$ws.Cells.Item(23, 9).Value = 1.02740799860509
$ws.Cells.Item(23, 10).Value = 1.024017392148712
$ws.Cells.Item(23, 11).Value = 1.026219232796392
$ws.Cells.Item(23, 12).Value = 1.022274199419794
$ws.Cells.Item(23, 13).Value = 1.016403692948373
$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.019764052320646
$ws.Cells.Item(24, 4).Value = 1.024421338347303
$ws.Cells.Item(24, 5).Value = 1.020802029546396
$ws.Cells.Item(24, 6).Value = 1.016235996499964
$ws.Cells.Item(24, 9).Value = 1.028052788393474
$ws.Cells.Item(24, 10).Value = 1.026064829130785
$ws.Cells.Item(24, 11).Value = 1.027832309059089
$ws.Cells.Item(24, 12).Value = 1.024225997963028
$ws.Cells.Item(24, 13).Value = 1.019676508389528
$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.022746841975043
$ws.Cells.Item(25, 4).Value = 1.026602088188735
$ws.Cells.Item(25, 5).Value = 1.023379373133355
$ws.Cells.Item(25, 6).Value = 1.020345577786254
$ws.Cells.Item(25, 9).Value = 1.028784941805884
$ws.Cells.Item(25, 10).Value = 1.02841939385682
$ws.Cells.Item(25, 11).Value = 1.029683807264326
$ws.Cells.Item(25, 12).Value = 1.026471455458504
$ws.Cells.Item(25, 13).Value = 1.023447481378648
